# Applies the cryptos-list refresh described by the commit:
#   "Updated cryptos list on Thu Aug 29 04:40:50 UTC 2024 with GitHub Actions"
#
# Coin/Link/Price/Volume(1h) live in columns B/C/D/E of the one data sheet.
# Most rows only get a refreshed Price (D) and Volume (E); a couple of rows
# (29/30 and 34/35) swapped ranking order, so their Coin/Link/Price all change.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = "59.056.50"; E = "  -0.49%  " }
    @{ Row = 3; D = "2.498.85"; E = "  +1.65%  " }
    @{ Row = 4; D = "1.00"; E = "  +0.05%  " }
    @{ Row = 5; D = "539.65"; E = "  +0.05%  " }
    @{ Row = 6; D = "143.54"; E = "  -2.93%  " }
    @{ Row = 7; D = "0.999"; E = "  +0.19%  " }
    @{ Row = 8; E = "  +0.41%  " }
    @{ Row = 9; D = "2.521.14"; E = "  +1.76%  " }
    @{ Row = 10; D = "0.100"; E = "  +0.90%  " }
    @{ Row = 11; E = "  +0.67%  " }
    @{ Row = 12; D = "5.56"; E = "  +4.41%  " }
    @{ Row = 13; D = "0.354"; E = "  +0.27%  " }
    @{ Row = 14; D = "2.964.53"; E = "  +2.27%  " }
    @{ Row = 15; D = "23.37"; E = "  -3.26%  " }
    @{ Row = 16; D = "59.009.04"; E = "  -0.45%  " }
    @{ Row = 17; E = "  +0.91%  " }
    @{ Row = 18; D = "2.522.52"; E = "  +0.10%  " }
    @{ Row = 19; D = "11.18"; E = "  +0.10%  " }
    @{ Row = 20; D = "4.27"; E = "  -1.96%  " }
    @{ Row = 21; D = "323.74"; E = "  -0.24%  " }
    @{ Row = 22; D = "1.00"; E = "  +3.33%  " }
    @{ Row = 23; D = "5.75"; E = "  +0.19%  " }
    @{ Row = 24; D = "61.91"; E = "  +2.16%  " }
    @{ Row = 25; D = "0.438"; E = "  -5.08%  " }
    @{ Row = 26; E = "  +0.58%  " }
    @{ Row = 27; D = "0.997"; E = "  +1.78%  " }
    @{ Row = 28; D = "7.84"; E = "  +1.70%  " }
    @{ Row = 29; B = "PEPE"; C = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"; D = "0.0₃0771"; E = "  -0.30%  " }
    @{ Row = 30; B = "PancakeSwap"; C = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"; D = "1.81"; E = "  -1.03%  " }
    @{ Row = 31; D = "6.61"; E = "  -2.39%  " }
    @{ Row = 32; E = "  -7.74%  " }
    @{ Row = 33; D = "0.999"; E = "  +0.14%  " }
    @{ Row = 34; B = "ImmutableX"; C = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"; D = "1.45"; E = "  +5.63%  " }
    @{ Row = 35; B = "Monero"; C = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"; D = "157.11"; E = "  -0.93%  " }
    @{ Row = 36; D = "18.59"; E = "  +1.08%  " }
    @{ Row = 37; D = "4.32"; E = "  -5.66%  " }
    @{ Row = 38; E = "  -9.88%  " }
    @{ Row = 39; D = "5.65"; E = "  -5.37%  " }
    @{ Row = 40; D = "36.84"; E = "  +0.19%  " }
    @{ Row = 41; D = "295.33"; E = "  -8.02%  " }
    @{ Row = 42; D = "3.67" }
    @{ Row = 43; E = "  -2.84%  " }
    @{ Row = 44; D = "0.995"; E = "  -0.03%  " }
    @{ Row = 45; D = "0.597"; E = "  +1.93%  " }
    @{ Row = 46; E = "  +0.40%  " }
    @{ Row = 47; D = "0.0927"; E = "  -1.47%  " }
    @{ Row = 48; E = "  +0.59%  " }
    @{ Row = 49; D = "18.49"; E = "  -0.45%  " }
    @{ Row = 50; E = "  -0.93%  " }
    @{ Row = 51; E = "  -2.61%  " }
)

foreach ($u in $updates) {
    if ($u.ContainsKey("B")) { $ws.Range("B" + $u.Row).Value = $u.B }
    if ($u.ContainsKey("C")) { $ws.Range("C" + $u.Row).Value = $u.C }
    if ($u.ContainsKey("D")) {
        # Price column holds text, not numbers (e.g. "1.00", "59.056.50").
        # A plain .Value assignment lets Excel auto-coerce a numeric-looking
        # string into a real number (dropping the trailing zero, mangling the
        # thousands-dot formatting, introducing float noise, ...). Marking the
        # cell as Text first keeps the write a literal string; resetting the
        # style back to Normal afterwards avoids leaving a stray number format
        # behind on the cell.
        $cell = $ws.Range("D" + $u.Row)
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
        $cell.Style = "Normal"
    }
    if ($u.ContainsKey("E")) { $ws.Range("E" + $u.Row).Value = $u.E }
}
